# Automatische test-sync: 2025-08-18 20:56:50
$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# Append new row 8 to the "Logs" sheet with the new mail-log entry.
$logs.Range("A8").Value = "Geen onderwerp"
$logs.Range("B8").Value = "onbekend"
$logs.Range("D8").Value = "Onbekend"
$logs.Range("F8").Value = "2025-08-18 20:56:03"
$logs.Range("G8").Value = "Nee"
$logs.Range("H8").Value = "Ja"
$logs.Range("I8").Value = "Nee"
$logs.Range("J8").Value = "Nee"

# Update the "Onbekend" count on the Dashboard sheet (3 -> 4).
$dashboard.Range("B2").Value = 4

# Extend the conditional-formatting ranges on "Logs" so the new row is covered.
$logs.Range("D2:D7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D8"))
$logs.Range("G2:G7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G8"))
$logs.Range("H2:H7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H8"))
$logs.Range("I2:I7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I8"))
$logs.Range("J2:J7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J8"))
